# Adding MPA test automation upload file
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Main Asset Number (ANLN1, col C) updates
$ws.Range("C6").Value = 60000393
$ws.Range("C7").Value = 60000393
$ws.Range("C10").Value = 60000393
$ws.Range("C15").Value = 60000394

# Asset Subnumber (ANLN2, col D) updates
$ws.Range("D8").Value = 280
$ws.Range("D9").Value = 280
$ws.Range("D11").Value = 280
$ws.Range("D12").Value = 280
$ws.Range("D13").Value = 280
$ws.Range("D14").Value = 280
$ws.Range("D16").Value = 281

# Row 15: add Revenue from asset sale (BF_ERLBT, col O)
$ws.Range("O15").Value = 10

# Row 16: Asset Transaction Type (BWASL, col E) changes from "260" (text/lookup) to numeric 250
$ws.Range("E16").Value = 250
